# Update cryptocurrency price/volume snapshot values (GitHub Actions data refresh).
# D/E-column numeric-looking text (e.g. "1.00", "214.53") must stay TEXT, matching the
# workbook's original inlineStr cells, so a leading apostrophe forces text entry for those.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "90.664.94"
$ws.Range("E2").Value = "  +1.32%  "

$ws.Range("D3").Value = "3.147.70"
$ws.Range("E3").Value = "  +3.81%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").Value = "'214.53"
$ws.Range("E5").Value = "  +2.11%  "

$ws.Range("D6").Value = "'624.46"
$ws.Range("E6").Value = "  +1.99%  "

$ws.Range("D7").Value = "'1.15"
$ws.Range("E7").Value = "  +30.35%  "

$ws.Range("D8").Value = "'0.368"
$ws.Range("E8").Value = "  +2.31%  "

$ws.Range("D10").Value = "3.143.50"
$ws.Range("E10").Value = "  +3.75%  "

$ws.Range("D11").Value = "'0.763"
$ws.Range("E11").Value = "  +15.03%  "

$ws.Range("D12").Value = "'0.202"
$ws.Range("E12").Value = "  +8.05%  "

$ws.Range("D13").Value = "'5.68"
$ws.Range("E13").Value = "  +6.31%  "

$ws.Range("D14").Value = "'0.0000244"
$ws.Range("E14").Value = "  +3.35%  "

$ws.Range("D15").Value = "'34.97"
$ws.Range("E15").Value = "  +9.84%  "

$ws.Range("D16").Value = "90.342.86"
$ws.Range("E16").Value = "  +2.05%  "

$ws.Range("D17").Value = "3.722.19"
$ws.Range("E17").Value = "  +3.66%  "

$ws.Range("D18").Value = "3.128.32"
$ws.Range("E18").Value = "  +3.48%  "

$ws.Range("D19").Value = "'3.67"
$ws.Range("E19").Value = "  +9.84%  "

$ws.Range("D20").Value = "'14.23"
$ws.Range("E20").Value = "  +6.71%  "

$ws.Range("D21").Value = "'462.05"
$ws.Range("E21").Value = "  +8.68%  "

$ws.Range("D22").Value = "'0.0000210"
$ws.Range("E22").Value = "  -1.54%  "

$ws.Range("D23").Value = "'9.04"
$ws.Range("E23").Value = "  +10.96%  "

$ws.Range("E24").Value = "  +6.07%  "

$ws.Range("D25").Value = "'5.86"
$ws.Range("E25").Value = "  +9.28%  "

$ws.Range("D26").Value = "'89.81"
$ws.Range("E26").Value = "  +7.80%  "

$ws.Range("D27").Value = "'12.03"
$ws.Range("E27").Value = "  +3.29%  "

$ws.Range("D28").Value = "3.314.56"
$ws.Range("E28").Value = "  +3.97%  "

$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("E30").Value = "  +1.47%  "

$ws.Range("E31").Value = "  -0.48%  "

$ws.Range("D32").Value = "'9.13"
$ws.Range("E32").Value = "  +11.63%  "

$ws.Range("D33").Value = "'27.22"
$ws.Range("E33").Value = "  +19.72%  "

$ws.Range("D34").Value = "'516.07"
$ws.Range("E34").Value = "  +2.97%  "

$ws.Range("E35").Value = "  +35.43%  "

$ws.Range("B36").Value = "PancakeSwap"
$ws.Range("C36").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D36").Value = "'1.92"
$ws.Range("E36").Value = "  +7.34%  "

$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "'3.61"
$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("D38").Value = "'0.141"
$ws.Range("E38").Value = "  +8.22%  "

$ws.Range("D39").Value = "'6.85"
$ws.Range("E39").Value = "  +3.60%  "

$ws.Range("E40").Value = "  +4.75%  "

$ws.Range("D41").Value = "'0.0867"
$ws.Range("E41").Value = "  +28.59%  "

$ws.Range("D42").Value = "'22.20"
$ws.Range("E42").Value = "  -0.14%  "

$ws.Range("E43").Value = "  +0.07%  "

$ws.Range("D44").Value = "'0.411"
$ws.Range("E44").Value = "  +13.99%  "

$ws.Range("D45").Value = "'1.95"
$ws.Range("E45").Value = "  +7.61%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "'148.84"
$ws.Range("E47").Value = "  +2.14%  "

$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").Value = "'4.58"
$ws.Range("E48").Value = "  +13.31%  "

$ws.Range("D49").Value = "'45.30"
$ws.Range("E49").Value = "  +4.68%  "

$ws.Range("E50").Value = "  +11.13%  "

$ws.Range("E51").Value = "  +13.94%  "
